$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The updated results table adds a 9th model block (LogisticRegression,
# rows 42-46), so first merge its column-A cell range the same way the
# other per-model blocks are merged, then clone the header-row-2 cell
# formatting (bold, centered, boxed border, style index 1) onto the new
# rows 42-46 for columns A and B.
$ws.Range("A42:A46").Merge()
$ws.Range("A2:B6").Copy()
$ws.Range("A42:B46").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rewrite rows 2-46 with the updated "Prior Concussion" results
# (new model ordering + values), including the new LogisticRegression
# block and its accuracy / macro avg / weighted avg rows.
$ws.Cells.Item(2,1).Value = "LGBMClassifier"
$ws.Cells.Item(2,2).Value = "No Prior Concussion: 0"
$ws.Cells.Item(2,3).Value = 0.8685362969196534
$ws.Cells.Item(2,4).Value = 0.8530612244897959
$ws.Cells.Item(2,5).Value = 0.8596347786884143
$ws.Cells.Item(2,6).Value = 49
$ws.Cells.Item(3,1).Value = ""
$ws.Cells.Item(3,2).Value = "Prior Concussion: 1"
$ws.Cells.Item(3,3).Value = 0.846838213546288
$ws.Cells.Item(3,4).Value = 0.859090909090909
$ws.Cells.Item(3,5).Value = 0.8516318002896168
$ws.Cells.Item(3,6).Value = 44
$ws.Cells.Item(4,1).Value = ""
$ws.Cells.Item(4,2).Value = "accuracy"
$ws.Cells.Item(4,3).Value = 0.8559139784946236
$ws.Cells.Item(4,4).Value = 0.8559139784946236
$ws.Cells.Item(4,5).Value = 0.8559139784946236
$ws.Cells.Item(4,6).Value = 0.8559139784946236
$ws.Cells.Item(5,1).Value = ""
$ws.Cells.Item(5,2).Value = "macro avg"
$ws.Cells.Item(5,3).Value = 0.8576872552329707
$ws.Cells.Item(5,4).Value = 0.8560760667903524
$ws.Cells.Item(5,5).Value = 0.8556332894890156
$ws.Cells.Item(5,6).Value = 93
$ws.Cells.Item(6,1).Value = ""
$ws.Cells.Item(6,2).Value = "weighted avg"
$ws.Cells.Item(6,3).Value = 0.8582705370440827
$ws.Cells.Item(6,4).Value = 0.8559139784946236
$ws.Cells.Item(6,5).Value = 0.8558484233169402
$ws.Cells.Item(6,6).Value = 93
$ws.Cells.Item(7,1).Value = "RandomForestClassifier"
$ws.Cells.Item(7,2).Value = "No Prior Concussion: 0"
$ws.Cells.Item(7,3).Value = 0.847704991087344
$ws.Cells.Item(7,4).Value = 0.8775510204081632
$ws.Cells.Item(7,5).Value = 0.8619016653449643
$ws.Cells.Item(7,6).Value = 49
$ws.Cells.Item(8,1).Value = ""
$ws.Cells.Item(8,2).Value = "Prior Concussion: 1"
$ws.Cells.Item(8,3).Value = 0.858312447786132
$ws.Cells.Item(8,4).Value = 0.8227272727272726
$ws.Cells.Item(8,5).Value = 0.8394360990905441
$ws.Cells.Item(8,6).Value = 44
$ws.Cells.Item(9,1).Value = ""
$ws.Cells.Item(9,2).Value = "accuracy"
$ws.Cells.Item(9,3).Value = 0.8516129032258064
$ws.Cells.Item(9,4).Value = 0.8516129032258064
$ws.Cells.Item(9,5).Value = 0.8516129032258064
$ws.Cells.Item(9,6).Value = 0.8516129032258064
$ws.Cells.Item(10,1).Value = ""
$ws.Cells.Item(10,2).Value = "macro avg"
$ws.Cells.Item(10,3).Value = 0.853008719436738
$ws.Cells.Item(10,4).Value = 0.850139146567718
$ws.Cells.Item(10,5).Value = 0.8506688822177543
$ws.Cells.Item(10,6).Value = 93
$ws.Cells.Item(11,1).Value = ""
$ws.Cells.Item(11,2).Value = "weighted avg"
$ws.Cells.Item(11,3).Value = 0.8527235727512867
$ws.Cells.Item(11,4).Value = 0.8516129032258064
$ws.Cells.Item(11,5).Value = 0.8512727952891096
$ws.Cells.Item(11,6).Value = 93
$ws.Cells.Item(12,1).Value = "NeuralNetClassifier"
$ws.Cells.Item(12,2).Value = "No Prior Concussion: 0"
$ws.Cells.Item(12,3).Value = 0.6874934498871716
$ws.Cells.Item(12,4).Value = 0.8653061224489796
$ws.Cells.Item(12,5).Value = 0.7613981552700898
$ws.Cells.Item(12,6).Value = 49
$ws.Cells.Item(13,1).Value = ""
$ws.Cells.Item(13,2).Value = "Prior Concussion: 1"
$ws.Cells.Item(13,3).Value = 0.796369239622613
$ws.Cells.Item(13,4).Value = 0.5454545454545454
$ws.Cells.Item(13,5).Value = 0.6331356371938563
$ws.Cells.Item(13,6).Value = 44
$ws.Cells.Item(14,1).Value = ""
$ws.Cells.Item(14,2).Value = "accuracy"
$ws.Cells.Item(14,3).Value = 0.7139784946236559
$ws.Cells.Item(14,4).Value = 0.7139784946236559
$ws.Cells.Item(14,5).Value = 0.7139784946236559
$ws.Cells.Item(14,6).Value = 0.7139784946236559
$ws.Cells.Item(15,1).Value = ""
$ws.Cells.Item(15,2).Value = "macro avg"
$ws.Cells.Item(15,3).Value = 0.7419313447548923
$ws.Cells.Item(15,4).Value = 0.7053803339517625
$ws.Cells.Item(15,5).Value = 0.697266896231973
$ws.Cells.Item(15,6).Value = 93
$ws.Cells.Item(16,1).Value = ""
$ws.Cells.Item(16,2).Value = "weighted avg"
$ws.Cells.Item(16,3).Value = 0.739004576213617
$ws.Cells.Item(16,4).Value = 0.7139784946236559
$ws.Cells.Item(16,5).Value = 0.70071481338456
$ws.Cells.Item(16,6).Value = 93
$ws.Cells.Item(17,1).Value = "DecisionTreeClassifier"
$ws.Cells.Item(17,2).Value = "No Prior Concussion: 0"
$ws.Cells.Item(17,3).Value = 0.6882132813905231
$ws.Cells.Item(17,4).Value = 0.7755102040816326
$ws.Cells.Item(17,5).Value = 0.7226024384647601
$ws.Cells.Item(17,6).Value = 49
$ws.Cells.Item(18,1).Value = ""
$ws.Cells.Item(18,2).Value = "Prior Concussion: 1"
$ws.Cells.Item(18,3).Value = 0.7197823860727086
$ws.Cells.Item(18,4).Value = 0.5954545454545455
$ws.Cells.Item(18,5).Value = 0.6372805332194786
$ws.Cells.Item(18,6).Value = 44
$ws.Cells.Item(19,1).Value = ""
$ws.Cells.Item(19,2).Value = "accuracy"
$ws.Cells.Item(19,3).Value = 0.6903225806451613
$ws.Cells.Item(19,4).Value = 0.6903225806451613
$ws.Cells.Item(19,5).Value = 0.6903225806451613
$ws.Cells.Item(19,6).Value = 0.6903225806451613
$ws.Cells.Item(20,1).Value = ""
$ws.Cells.Item(20,2).Value = "macro avg"
$ws.Cells.Item(20,3).Value = 0.7039978337316158
$ws.Cells.Item(20,4).Value = 0.685482374768089
$ws.Cells.Item(20,5).Value = 0.6799414858421194
$ws.Cells.Item(20,6).Value = 93
$ws.Cells.Item(21,1).Value = ""
$ws.Cells.Item(21,2).Value = "weighted avg"
$ws.Cells.Item(21,3).Value = 0.7031492018853206
$ws.Cells.Item(21,4).Value = 0.6903225806451613
$ws.Cells.Item(21,5).Value = 0.6822350854454873
$ws.Cells.Item(21,6).Value = 93
$ws.Cells.Item(22,1).Value = "LinearBoostClassifier"
$ws.Cells.Item(22,2).Value = "No Prior Concussion: 0"
$ws.Cells.Item(22,3).Value = 0.6904106280193237
$ws.Cells.Item(22,4).Value = 0.6571428571428571
$ws.Cells.Item(22,5).Value = 0.6732510592119694
$ws.Cells.Item(22,6).Value = 49
$ws.Cells.Item(23,1).Value = ""
$ws.Cells.Item(23,2).Value = "Prior Concussion: 1"
$ws.Cells.Item(23,3).Value = 0.6387647754137116
$ws.Cells.Item(23,4).Value = 0.6727272727272727
$ws.Cells.Item(23,5).Value = 0.6551925358471523
$ws.Cells.Item(23,6).Value = 44
$ws.Cells.Item(24,1).Value = ""
$ws.Cells.Item(24,2).Value = "accuracy"
$ws.Cells.Item(24,3).Value = 0.664516129032258
$ws.Cells.Item(24,4).Value = 0.664516129032258
$ws.Cells.Item(24,5).Value = 0.664516129032258
$ws.Cells.Item(24,6).Value = 0.664516129032258
$ws.Cells.Item(25,1).Value = ""
$ws.Cells.Item(25,2).Value = "macro avg"
$ws.Cells.Item(25,3).Value = 0.6645877017165176
$ws.Cells.Item(25,4).Value = 0.6649350649350649
$ws.Cells.Item(25,5).Value = 0.6642217975295608
$ws.Cells.Item(25,6).Value = 93
$ws.Cells.Item(26,1).Value = ""
$ws.Cells.Item(26,2).Value = "weighted avg"
$ws.Cells.Item(26,3).Value = 0.6659760310876363
$ws.Cells.Item(26,4).Value = 0.664516129032258
$ws.Cells.Item(26,5).Value = 0.6647072417060345
$ws.Cells.Item(26,6).Value = 93
$ws.Cells.Item(27,1).Value = "XGBClassifier"
$ws.Cells.Item(27,2).Value = "No Prior Concussion: 0"
$ws.Cells.Item(27,3).Value = 0.8
$ws.Cells.Item(27,4).Value = 0.4122448979591836
$ws.Cells.Item(27,5).Value = 0.5395279673801024
$ws.Cells.Item(27,6).Value = 49
$ws.Cells.Item(28,1).Value = ""
$ws.Cells.Item(28,2).Value = "Prior Concussion: 1"
$ws.Cells.Item(28,3).Value = 0.6171252133400451
$ws.Cells.Item(28,4).Value = 1
$ws.Cells.Item(28,5).Value = 0.7597273097622192
$ws.Cells.Item(28,6).Value = 44
$ws.Cells.Item(29,1).Value = ""
$ws.Cells.Item(29,2).Value = "accuracy"
$ws.Cells.Item(29,3).Value = 0.6903225806451613
$ws.Cells.Item(29,4).Value = 0.6903225806451613
$ws.Cells.Item(29,5).Value = 0.6903225806451613
$ws.Cells.Item(29,6).Value = 0.6903225806451613
$ws.Cells.Item(30,1).Value = ""
$ws.Cells.Item(30,2).Value = "macro avg"
$ws.Cells.Item(30,3).Value = 0.7085626066700226
$ws.Cells.Item(30,4).Value = 0.7061224489795919
$ws.Cells.Item(30,5).Value = 0.6496276385711608
$ws.Cells.Item(30,6).Value = 93
$ws.Cells.Item(31,1).Value = ""
$ws.Cells.Item(31,2).Value = "weighted avg"
$ws.Cells.Item(31,3).Value = 0.713478595558731
$ws.Cells.Item(31,4).Value = 0.6903225806451613
$ws.Cells.Item(31,5).Value = 0.6437083014103513
$ws.Cells.Item(31,6).Value = 93
$ws.Cells.Item(32,1).Value = "ElasticNet"
$ws.Cells.Item(32,2).Value = "No Prior Concussion: 0"
$ws.Cells.Item(32,3).Value = 0.6430986863002752
$ws.Cells.Item(32,4).Value = 0.7061224489795919
$ws.Cells.Item(32,5).Value = 0.6728246923299885
$ws.Cells.Item(32,6).Value = 49
$ws.Cells.Item(33,1).Value = ""
$ws.Cells.Item(33,2).Value = "Prior Concussion: 1"
$ws.Cells.Item(33,3).Value = 0.6336562784420801
$ws.Cells.Item(33,4).Value = 0.5636363636363637
$ws.Cells.Item(33,5).Value = 0.5960855750286982
$ws.Cells.Item(33,6).Value = 44
$ws.Cells.Item(34,1).Value = ""
$ws.Cells.Item(34,2).Value = "accuracy"
$ws.Cells.Item(34,3).Value = 0.6387096774193548
$ws.Cells.Item(34,4).Value = 0.6387096774193548
$ws.Cells.Item(34,5).Value = 0.6387096774193548
$ws.Cells.Item(34,6).Value = 0.6387096774193548
$ws.Cells.Item(35,1).Value = ""
$ws.Cells.Item(35,2).Value = "macro avg"
$ws.Cells.Item(35,3).Value = 0.6383774823711776
$ws.Cells.Item(35,4).Value = 0.6348794063079778
$ws.Cells.Item(35,5).Value = 0.6344551336793434
$ws.Cells.Item(35,6).Value = 93
$ws.Cells.Item(36,1).Value = ""
$ws.Cells.Item(36,2).Value = "weighted avg"
$ws.Cells.Item(36,3).Value = 0.6386313105394088
$ws.Cells.Item(36,4).Value = 0.6387096774193548
$ws.Cells.Item(36,5).Value = 0.6365180131766899
$ws.Cells.Item(36,6).Value = 93
$ws.Cells.Item(37,1).Value = "SVC"
$ws.Cells.Item(37,2).Value = "No Prior Concussion: 0"
$ws.Cells.Item(37,3).Value = 0.5934526856559148
$ws.Cells.Item(37,4).Value = 0.6612244897959184
$ws.Cells.Item(37,5).Value = 0.6232929866098293
$ws.Cells.Item(37,6).Value = 49
$ws.Cells.Item(38,1).Value = ""
$ws.Cells.Item(38,2).Value = "Prior Concussion: 1"
$ws.Cells.Item(38,3).Value = 0.59195388164417
$ws.Cells.Item(38,4).Value = 0.509090909090909
$ws.Cells.Item(38,5).Value = 0.5434946657915315
$ws.Cells.Item(38,6).Value = 44
$ws.Cells.Item(39,1).Value = ""
$ws.Cells.Item(39,2).Value = "accuracy"
$ws.Cells.Item(39,3).Value = 0.589247311827957
$ws.Cells.Item(39,4).Value = 0.589247311827957
$ws.Cells.Item(39,5).Value = 0.589247311827957
$ws.Cells.Item(39,6).Value = 0.589247311827957
$ws.Cells.Item(40,1).Value = ""
$ws.Cells.Item(40,2).Value = "macro avg"
$ws.Cells.Item(40,3).Value = 0.5927032836500424
$ws.Cells.Item(40,4).Value = 0.5851576994434138
$ws.Cells.Item(40,5).Value = 0.5833938262006805
$ws.Cells.Item(40,6).Value = 93
$ws.Cells.Item(41,1).Value = ""
$ws.Cells.Item(41,2).Value = "weighted avg"
$ws.Cells.Item(41,3).Value = 0.5927435740804656
$ws.Cells.Item(41,4).Value = 0.589247311827957
$ws.Cells.Item(41,5).Value = 0.5855389423517099
$ws.Cells.Item(41,6).Value = 93
$ws.Cells.Item(42,1).Value = "LogisticRegression"
$ws.Cells.Item(42,2).Value = "No Prior Concussion: 0"
$ws.Cells.Item(42,3).Value = 0.4915584415584416
$ws.Cells.Item(42,4).Value = 0.5102040816326531
$ws.Cells.Item(42,5).Value = 0.499760612005901
$ws.Cells.Item(42,6).Value = 49
$ws.Cells.Item(43,1).Value = ""
$ws.Cells.Item(43,2).Value = "Prior Concussion: 1"
$ws.Cells.Item(43,3).Value = 0.5560418685517033
$ws.Cells.Item(43,4).Value = 0.6409090909090909
$ws.Cells.Item(43,5).Value = 0.5782509282044364
$ws.Cells.Item(43,6).Value = 44
$ws.Cells.Item(44,1).Value = ""
$ws.Cells.Item(44,2).Value = "accuracy"
$ws.Cells.Item(44,3).Value = 0.5720430107526882
$ws.Cells.Item(44,4).Value = 0.5720430107526882
$ws.Cells.Item(44,5).Value = 0.5720430107526882
$ws.Cells.Item(44,6).Value = 0.5720430107526882
$ws.Cells.Item(45,1).Value = ""
$ws.Cells.Item(45,2).Value = "macro avg"
$ws.Cells.Item(45,3).Value = 0.5238001550550725
$ws.Cells.Item(45,4).Value = 0.5755565862708719
$ws.Cells.Item(45,5).Value = 0.5390057701051687
$ws.Cells.Item(45,6).Value = 93
$ws.Cells.Item(46,1).Value = ""
$ws.Cells.Item(46,2).Value = "weighted avg"
$ws.Cells.Item(46,3).Value = 0.5220667295982644
$ws.Cells.Item(46,4).Value = 0.5720430107526882
$ws.Cells.Item(46,5).Value = 0.5368958153686488
$ws.Cells.Item(46,6).Value = 93
